# Rename the three logo inline pictures in the document's header/footer
# parts (swap image1.png <-> image2.png for the two Pearson logos that are
# currently both called "image1.png", and image2.jpg -> image1.jpg for the
# BTEC logo), matching the authoring tool's re-numbering of the media parts.
#
# wdHeaderFooterIndex constants used below:
#   1 = wdHeaderFooterPrimary   (the document's main/default header or footer)
#   2 = wdHeaderFooterFirstPage (the first-page header or footer)

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- First-page footer: Pearson logo "image1.png" -> "image2.png" ---
$footerFirst = $sec.Footers.Item(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -gt 0) {
    $shp = $footerFirst.Range.InlineShapes.Item(1)
    $shp.Name = "image2.png"
}

# --- Default (primary) footer: Pearson logo "image1.png" -> "image2.png" ---
$footerPrimary = $sec.Footers.Item(1)
if ($footerPrimary.Exists -and $footerPrimary.Range.InlineShapes.Count -gt 0) {
    $shp = $footerPrimary.Range.InlineShapes.Item(1)
    $shp.Name = "image2.png"
}

# --- First-page header: BTEC logo "image2.jpg" -> "image1.jpg" ---
$headerFirst = $sec.Headers.Item(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -gt 0) {
    $shp = $headerFirst.Range.InlineShapes.Item(1)
    $shp.Name = "image1.jpg"
}
